$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto list snapshot. Values are written as literal text (matching the
# workbook's existing text cells, e.g. dotted "thousands" price strings
# like "1.672.16" and padded percentage strings like "  -3.72%  ").
#
# Excel's Range.Value setter auto-detects numeric-looking strings and
# coerces them to numbers (which would silently drop significant trailing
# zeros, e.g. "0.06370" -> 0.0637). Toggling NumberFormat to Text ("@")
# for the write forces the literal string to be kept, and restoring the
# cell's style back to "Normal" afterwards leaves formatting untouched.
$cellValues = [ordered]@{
    "D2" = "26.331.95"
    "E2" = "  -5.53%  "
    "D3" = "1.672.16"
    "E3" = "  -3.72%  "
    "E4" = "  +0.23%  "
    "D5" = "217.87"
    "E5" = "  -3.74%  "
    "D6" = "0.5092"
    "E6" = "  -11.63%  "
    "E7" = "  +0.17%  "
    "D8" = "0.2663"
    "E8" = "  -2.59%  "
    "D9" = "0.06370"
    "E9" = "  -3.64%  "
    "D10" = "21.57"
    "E10" = "  -6.35%  "
    "D11" = "0.07362"
    "E11" = "  -2.40%  "
    "D12" = "1.672.35"
    "E12" = "  -3.70%  "
    "D13" = "4.556"
    "E13" = "  -3.08%  "
    "D14" = "0.5829"
    "E14" = "  -3.05%  "
    "D15" = "1.897.68"
    "E15" = "  -3.82%  "
    "D16" = "0.000008545"
    "E16" = "  -1.30%  "
    "D17" = "64.89"
    "E17" = "  -12.78%  "
    "D18" = "26.385.04"
    "E18" = "  -5.29%  "
    "D19" = "4.951"
    "E19" = "  -6.76%  "
    "E20" = "  +0.25%  "
    "D21" = "10.86"
    "E21" = "  -3.57%  "
    "D22" = "189.48"
    "E22" = "  -7.54%  "
    "D23" = "6.216"
    "E23" = "  -5.91%  "
    "D24" = "1.005"
    "E24" = "  +0.19%  "
    "D25" = "143.44"
    "E25" = "  -4.41%  "
    "D26" = "7.693"
    "E26" = "  -4.02%  "
    "D27" = "0.1181"
    "E27" = "  -4.09%  "
    "D28" = "15.66"
    "E28" = "  -3.17%  "
    "D29" = "0.05890"
    "E29" = "  -4.52%  "
    "D30" = "1.262"
    "E30" = "  -8.72%  "
    "D31" = "1.323"
    "E31" = "  -4.88%  "
    "D32" = "3.533"
    "E32" = "  -5.17%  "
    "D33" = "3.518"
    "E33" = "  -5.85%  "
    "D34" = "1.644"
    "E34" = "  -1.76%  "
    "D35" = "1.014"
    "E35" = "  -1.81%  "
    "D36" = "0.6013"
    "E36" = "  -6.51%  "
    "D37" = "2.354"
    "E37" = "  -3.22%  "
    "D38" = "2.649"
    "E38" = "  -2.12%  "
    "D39" = "0.01618"
    "E39" = "  -2.95%  "
    "D40" = "6.021"
    "E40" = "  -2.33%  "
    "D41" = "1.077.20"
    "E41" = "  -3.79%  "
    "D42" = "0.8658"
    "E42" = "  -0.73%  "
    "E43" = "  +0.48%  "
    "D44" = "99.75"
    "E44" = "  -0.26%  "
    "D45" = "1.820.47"
    "E45" = "  -3.48%  "
    "D46" = "0.00000000111"
    "E46" = "  +3.19%  "
    "D47" = "55.95"
    "E47" = "  -5.63%  "
    "D48" = "1.005"
    "E48" = "  +0.27%  "
    "D49" = "8.092"
    "E49" = "  -1.67%  "
    "D51" = "0.05184"
    "E51" = "  -3.55%  "
}

foreach ($addr in $cellValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$addr]
    $cell.Style = "Normal"
}
